# Applies the "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# edit: replaces the account-statement detail rows with a new data set
# (new employee added, periods re-ordered/changed) and updates the
# summary header figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room: the detail table grows from 12 data rows (16-27) to
#    14 data rows (16-29). Insert two blank rows just above the old
#    last (heavy-bottom-border) row so that row shifts down intact
#    with its special formatting, and copy the "normal" row format
#    into the two newly inserted rows.
# ------------------------------------------------------------------
$ws.Rows("27:28").Insert(-4121)   # xlShiftDown

$ws.Range("B26:J26").Copy()
$ws.Range("B27:J28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Summary / header fields
# ------------------------------------------------------------------
$ws.Range("E11").Value = 733653
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 13

# ------------------------------------------------------------------
# 3. Detail table body (rows 16-29)
# ------------------------------------------------------------------
$rows = @(
    @{r=16; b="CC"; c="1047408226"; d="MARIA ANGELICA RODRIGUEZ ORTIZ";        e="2203"; f=60000; g=1500000},
    @{r=17; b="CC"; c="1047408226"; d="MARIA ANGELICA RODRIGUEZ ORTIZ";        e="2202"; f=60000; g=1500000},
    @{r=18; b="CC"; c="1047408226"; d="MARIA ANGELICA RODRIGUEZ ORTIZ";        e="2201"; f=52000; g=1500000},
    @{r=19; b="CC"; c="1143337876"; d="BEATRIZ YADIRA RUBIO PERALTA";          e="2506"; f=50400; g=600000},
    @{r=20; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1902"; f=48000; g=1200000},
    @{r=21; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1901"; f=48000; g=1200000},
    @{r=22; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1812"; f=48000; g=1200000},
    @{r=23; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1811"; f=48000; g=1200000},
    @{r=24; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1810"; f=48000; g=1200000},
    @{r=25; b="CC"; c="73009628";   d="JOSE LUIS LUNA OVIEDO";                 e="1809"; f=48000; g=1200000},
    @{r=26; b="CC"; c="1047481580"; d="MARYELIS BEATRIZ MORA DE LA ESPRIELLA"; e="2507"; f=59800; g=1495000},
    @{r=27; b="CC"; c="1047481580"; d="MARYELIS BEATRIZ MORA DE LA ESPRIELLA"; e="2506"; f=59800; g=1495000},
    @{r=28; b="CC"; c="1047481580"; d="MARYELIS BEATRIZ MORA DE LA ESPRIELLA"; e="2505"; f=59800; g=1495000},
    @{r=29; b="CC"; c="1047481580"; d="MARYELIS BEATRIZ MORA DE LA ESPRIELLA"; e="2504"; f=43853; g=1495000}
)

foreach ($row in $rows) {
    $n = $row.r
    $ws.Range("B$n").Value = $row.b
    $ws.Range("C$n").Value = $row.c
    $ws.Range("D$n").Value = $row.d
    $ws.Range("E$n").Value = $row.e
    $ws.Range("F$n").Value = $row.f
    $ws.Range("G$n").Value = $row.g
}

Write-Host "Edit applied"
